$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "330.39"
Set-TextValue "E2" "2.75%"

Set-TextValue "D3" "41.07"
Set-TextValue "E3" "4.27%"

Set-TextValue "D4" "5.699"
Set-TextValue "E4" "-3.63%"

Set-TextValue "D5" "0.08170"
Set-TextValue "E5" "1.94%"

Set-TextValue "D6" "2.057"
Set-TextValue "E6" "6.64%"

Set-TextValue "E7" "1.23%"

Set-TextValue "D8" "2.949"
Set-TextValue "E8" "0.23%"

Set-TextValue "D9" "0.9254"

Set-TextValue "D10" "0.1251"
Set-TextValue "E10" "-0.16%"

Set-TextValue "D11" "0.1956"
Set-TextValue "E11" "-0.17%"

Set-TextValue "D12" "0.09316"
Set-TextValue "E12" "1.68%"

Set-TextValue "D13" "0.03694"
Set-TextValue "E13" "4.00%"

Set-TextValue "E14" "10.12%"

Set-TextValue "D15" "0.001300"
Set-TextValue "E15" "0.18%"

Set-TextValue "D16" "0.006194"
Set-TextValue "E16" "-0.88%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D17" "3.433"
Set-TextValue "E17" "2.29%"

$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D18" "4.542"
Set-TextValue "E18" "-0.39%"

$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D19" "0.3485"
Set-TextValue "E19" "-1.46%"

$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D20" "8.343"
Set-TextValue "E20" "-4.97%"

$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D21" "0.1381"
Set-TextValue "E21" "-3.38%"

$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D22" "0.2653"
Set-TextValue "E22" "10.12%"

$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D23" "0.04454"
Set-TextValue "E23" "0.20%"

$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D24" "0.001270"
Set-TextValue "E24" "0.68%"

$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D25" "0.004291"
Set-TextValue "E25" "-2.33%"

Set-TextValue "D39" "0.02801"
Set-TextValue "E39" "16.52%"

Set-TextValue "D40" "0.05478"
Set-TextValue "E40" "5.85%"

Set-TextValue "D41" "0.007672"
Set-TextValue "E41" "3.16%"

Set-TextValue "D42" "0.009437"
Set-TextValue "E42" "0.99%"

Set-TextValue "D43" "0.1417"
Set-TextValue "E43" "0.84%"

Set-TextValue "D44" "0.002134"
Set-TextValue "E44" "0.61%"

Set-TextValue "D45" "0.01148"
Set-TextValue "E45" "2.15%"

Set-TextValue "D46" "0.00006864"
Set-TextValue "E46" "1.93%"

Set-TextValue "D47" "0.00000000751"
Set-TextValue "E47" "0.14%"

Set-TextValue "D49" "0.003530"
Set-TextValue "E49" "17.40%"

Set-TextValue "E50" "0.14%"

Set-TextValue "E51" "0.14%"
